$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new postal code labels in column B for rows 17, 19, 21, 23, 25
$ws.Range("B17").Value = "1. M8X 1E9"
$ws.Range("B19").Value = "2. M5A 2L2"
$ws.Range("B21").Value = "3. M6K 1L4"
$ws.Range("B23").Value = "4. 89109"
$ws.Range("B25").Value = "5. M6K"

# Update the view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("B23").Select()
